$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add "ulna" to B7 and "radius" to B8, matching the pattern of the existing
# single-word bone terms already present in column B (tibia, pes, femur, humerus)
$ws.Range("B7").Value = "ulna"
$ws.Range("B8").Value = "radius"

# Update the active cell selection to reflect the new last-used cell
$ws.Range("B8").Select()
